$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# (2026-02-28 -> 2026-03-01, serial 46081 -> 46082) for every data row
# (rows 2 through 447).
for ($r = 2; $r -le 447; $r++) {
    $ws.Cells.Item($r, 3).Value = 46082
}
